$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 4198
$ws.Range("B2").Value = "Cauã Viana"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Consulta medica"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45087
$ws.Range("G2").Value = 4406.91

# Row 3
$ws.Range("A3").Value = 16809
$ws.Range("B3").Value = "Felipe Costa"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Doenca"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 45091
$ws.Range("G3").Value = 9596.950000000001

# Row 4
$ws.Range("A4").Value = 92411
$ws.Range("B4").Value = "Dra. Maria Isis Almeida"
$ws.Range("C4").Value = "Juridico"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45084
$ws.Range("G4").Value = 3105.35

# Row 5
$ws.Range("A5").Value = 339
$ws.Range("B5").Value = "Alexandre Moura"
$ws.Range("D5").Value = "Viagem de negocios"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45086
$ws.Range("G5").Value = 2138.75

# Row 6
$ws.Range("A6").Value = 51864
$ws.Range("B6").Value = "Ana Júlia Mendes"
$ws.Range("C6").Value = "Financeiro"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 6
$ws.Range("F6").Value = 45080
$ws.Range("G6").Value = 8324.74

# Row 7
$ws.Range("A7").Value = 54210
$ws.Range("B7").Value = "Fernando Siqueira"
$ws.Range("C7").Value = "Operacoes"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 45095
$ws.Range("G7").Value = 9701.370000000001

# Row 8
$ws.Range("A8").Value = 5652
$ws.Range("B8").Value = "Arthur Miguel Teixeira"
$ws.Range("C8").Value = "Vendas"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45101
$ws.Range("G8").Value = 7767.19

# Row 9
$ws.Range("A9").Value = 61315
$ws.Range("B9").Value = "Lorena Campos"
$ws.Range("D9").Value = "Viagem de negocios"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45091
$ws.Range("G9").Value = 7960.95

# Row 10
$ws.Range("A10").Value = 47139
$ws.Range("B10").Value = "Aylla Pastor"
$ws.Range("C10").Value = "Financeiro"
$ws.Range("D10").Value = "Consulta medica"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45086
$ws.Range("G10").Value = 8247.07

# Row 11
$ws.Range("A11").Value = 70344
$ws.Range("B11").Value = "Sra. Ana Júlia Porto"
$ws.Range("C11").Value = "Marketing"
$ws.Range("D11").Value = "Doenca"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45086
$ws.Range("G11").Value = 3189.53
